$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting so
# numeric-looking values (e.g. "1.00", "382.40") are not coerced
# into numbers and lose their literal representation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "62.005.25"
$ws.Cells.Item(2, 5).Value = "  -1.75%  "

$ws.Cells.Item(3, 4).Value = "3.422.05"
$ws.Cells.Item(3, 5).Value = "  -1.06%  "

$ws.Cells.Item(5, 4).Value = "579.16"

$ws.Cells.Item(6, 4).Value = "154.02"
$ws.Cells.Item(6, 5).Value = "  +4.41%  "

$ws.Cells.Item(7, 5).Value = "  +0.03%  "

$ws.Cells.Item(8, 5).Value = "  +1.36%  "

$ws.Cells.Item(9, 5).Value = "  +2.58%  "

$ws.Cells.Item(10, 5).Value = "  +0.44%  "

$ws.Cells.Item(11, 4).Value = "0.419"
$ws.Cells.Item(11, 5).Value = "  +3.43%  "

$ws.Cells.Item(12, 4).Value = "4.010.85"
$ws.Cells.Item(12, 5).Value = "  -1.04%  "

$ws.Cells.Item(13, 5).Value = "  +0.62%  "

$ws.Cells.Item(14, 4).Value = "28.35"
$ws.Cells.Item(14, 5).Value = "  -2.91%  "

$ws.Cells.Item(15, 5).Value = "  +0.16%  "

$ws.Cells.Item(16, 4).Value = "3.417.13"
$ws.Cells.Item(16, 5).Value = "  -1.29%  "

$ws.Cells.Item(17, 4).Value = "62.032.69"
$ws.Cells.Item(17, 5).Value = "  -1.74%  "

$ws.Cells.Item(18, 4).Value = "6.59"
$ws.Cells.Item(18, 5).Value = "  +2.83%  "

$ws.Cells.Item(19, 4).Value = "14.47"
$ws.Cells.Item(19, 5).Value = "  +0.24%  "

$ws.Cells.Item(20, 4).Value = "8.94"
$ws.Cells.Item(20, 5).Value = "  -3.41%  "

$ws.Cells.Item(21, 4).Value = "382.40"
$ws.Cells.Item(21, 5).Value = "  -1.15%  "

$ws.Cells.Item(22, 4).Value = "0.572"
$ws.Cells.Item(22, 5).Value = "  +1.71%  "

$ws.Cells.Item(23, 4).Value = "75.81"
$ws.Cells.Item(23, 5).Value = "  +1.81%  "

$ws.Cells.Item(24, 5).Value = "  +0.01%  "

$ws.Cells.Item(25, 4).Value = "3.561.03"
$ws.Cells.Item(25, 5).Value = "  -1.16%  "

$ws.Cells.Item(26, 5).Value = "  -1.75%  "

$ws.Cells.Item(27, 5).Value = "  -1.20%  "

$ws.Cells.Item(28, 4).Value = "7.64"
$ws.Cells.Item(28, 5).Value = "  +0.09%  "

$ws.Cells.Item(29, 4).Value = "1.00"
$ws.Cells.Item(29, 5).Value = "  +0.05%  "

$ws.Cells.Item(30, 5).Value = "  +0.45%  "

$ws.Cells.Item(31, 4).Value = "7.88"
$ws.Cells.Item(31, 5).Value = "  -3.33%  "

$ws.Cells.Item(32, 5).Value = "  -0.05%  "

$ws.Cells.Item(33, 2).Value = "EthereumClassic"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(33, 4).Value = "23.28"
$ws.Cells.Item(33, 5).Value = "  -0.30%  "

$ws.Cells.Item(34, 2).Value = "Fetch.AI"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(34, 4).Value = "1.33"
$ws.Cells.Item(34, 5).Value = "  -0.43%  "

$ws.Cells.Item(35, 4).Value = "5.57"

$ws.Cells.Item(36, 5).Value = "  +0.69%  "

$ws.Cells.Item(37, 4).Value = "6.97"
$ws.Cells.Item(37, 5).Value = "  -2.29%  "

$ws.Cells.Item(38, 4).Value = "31.12"
$ws.Cells.Item(38, 5).Value = "  -2.26%  "

$ws.Cells.Item(39, 4).Value = "168.12"
$ws.Cells.Item(39, 5).Value = "  -0.11%  "

$ws.Cells.Item(40, 4).Value = "3.456.79"
$ws.Cells.Item(40, 5).Value = "  -1.13%  "

$ws.Cells.Item(41, 4).Value = "0.0786"
$ws.Cells.Item(41, 5).Value = "  +2.20%  "

$ws.Cells.Item(42, 4).Value = "42.71"
$ws.Cells.Item(42, 5).Value = "  +0.73%  "

$ws.Cells.Item(43, 4).Value = "0.781"
$ws.Cells.Item(43, 5).Value = "  -1.45%  "

$ws.Cells.Item(44, 4).Value = "4.43"
$ws.Cells.Item(44, 5).Value = "  +1.62%  "

$ws.Cells.Item(45, 4).Value = "1.68"
$ws.Cells.Item(45, 5).Value = "  -2.92%  "

$ws.Cells.Item(46, 4).Value = "1.16"
$ws.Cells.Item(46, 5).Value = "  -3.40%  "

$ws.Cells.Item(47, 4).Value = "2.551.98"
$ws.Cells.Item(47, 5).Value = "  -1.54%  "

$ws.Cells.Item(48, 4).Value = "6.82"
$ws.Cells.Item(48, 5).Value = "  +0.28%  "

$ws.Cells.Item(49, 4).Value = "23.08"
$ws.Cells.Item(49, 5).Value = "  +0.89%  "

$ws.Cells.Item(50, 4).Value = "2.18"
$ws.Cells.Item(50, 5).Value = "  -5.37%  "

$ws.Cells.Item(51, 5).Value = "  -0.01%  "
